$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.563.51"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.673.14"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'219.70"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'29.57"
$ws.Range("E8").Value = "  +3.87%  "
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("D10").Value = "'0.0644"
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.912.80"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.614"
$ws.Range("E13").Value = "  +9.13%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.26"
$ws.Range("E14").Value = "  +10.99%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.663.37"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "'4.01"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("D17").Value = "30.596.60"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "'66.41"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").Value = "'243.24"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  +3.35%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D23").Value = "'10.01"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'158.36"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "'15.88"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("E27").Value = "  +2.65%  "
$ws.Range("D28").Value = "'6.69"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  +3.15%  "
$ws.Range("D32").Value = "'3.47"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("D34").Value = "1.487.53"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").Value = "  +6.91%  "
$ws.Range("D36").Value = "'84.89"
$ws.Range("E36").Value = "  +12.13%  "
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "'0.600"
$ws.Range("E38").Value = "  +8.66%  "
$ws.Range("E39").Value = "  +5.26%  "
$ws.Range("E40").Value = "  -4.49%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "'0.840"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'51.48"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'5.51"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("D49").Value = "1.805.73"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("D50").Value = "'94.83"
$ws.Range("E50").Value = "  +4.32%  "
$ws.Range("E51").Value = "  -0.97%  "
